$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados..." timestamp refresh (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 23:10"

# --- Row 4: Estados Unidos - refreshed totals ---
$ws.Range("B4").Value = 1787514
$ws.Range("C4").Value = 19053
$ws.Range("D4").Value = 503933
$ws.Range("E4").Value = 1179259
$ws.Range("G4").Value = 992
$ws.Range("H4").Value = 104322

# --- Row 146: Ruanda - refreshed totals ---
$ws.Range("B146").Value = 355
$ws.Range("C146").Value = 6
$ws.Range("D146").Value = 247
$ws.Range("E146").Value = 108

# --- Rows 171/172: Libia overtakes Trinidad y Tobago in the ranking ---
$ws.Range("A171").Value = "Libia"
$ws.Range("B171").Value = 118
$ws.Range("C171").Value = 13
$ws.Range("D171").Value = 41
$ws.Range("E171").Value = 72
$ws.Range("H171").Value = 5

$ws.Range("A172").Value = "Trinidad yTobago"
$ws.Range("B172").Value = 116
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 108
$ws.Range("E172").Value = 0
$ws.Range("H172").Value = 8

# --- Rows 200/201: Santa Lucia overtakes Belice in the ranking ---
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# --- Rows 215/216: San Bartolome overtakes Bonaire, San Eustaquio y Saba ---
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
